$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Map of cell address -> new "Recorded By" value (names reordered)
$updates = @{
    "G2" = "system, System, backup@backdoor.com"
    "G4" = "System, backup@backdoor.com"
    "G5" = "System, backup@backdoor.com"
    "G7" = "System, admin@admin.com"
    "G8" = "System, backup@backdoor.com"
    "G11" = "dnasr281@gmail.com, System"
    "G17" = "dnasr281@gmail.com, System"
    "G28" = "system, System, backup@backdoor.com"
    "G30" = "System, backup@backdoor.com"
    "G31" = "System, backup@backdoor.com"
    "G33" = "System, admin@admin.com"
    "G34" = "System, backup@backdoor.com"
    "G37" = "dnasr281@gmail.com, System"
    "G43" = "dnasr281@gmail.com, System"
    "G54" = "system, System, backup@backdoor.com"
    "G56" = "System, backup@backdoor.com"
    "G57" = "System, backup@backdoor.com"
    "G59" = "System, admin@admin.com"
    "G60" = "System, backup@backdoor.com"
    "G63" = "dnasr281@gmail.com, System"
    "G69" = "dnasr281@gmail.com, System"
    "G80" = "System, backup@backdoor.com"
    "G81" = "System, backup@backdoor.com"
    "G82" = "System, backup@backdoor.com"
    "G87" = "dnasr281@gmail.com, admin@admin.com"
    "G93" = "dnasr281@gmail.com, System"
    "G94" = "dnasr281@gmail.com, System"
    "G96" = "dnasr281@gmail.com, System"
    "G106" = "System, backup@backdoor.com"
    "G107" = "System, backup@backdoor.com"
    "G108" = "System, backup@backdoor.com"
    "G113" = "dnasr281@gmail.com, admin@admin.com"
    "G119" = "dnasr281@gmail.com, System"
    "G120" = "dnasr281@gmail.com, System"
    "G122" = "dnasr281@gmail.com, System"
    "G132" = "System, backup@backdoor.com"
    "G133" = "System, backup@backdoor.com"
    "G134" = "System, backup@backdoor.com"
    "G139" = "dnasr281@gmail.com, admin@admin.com"
    "G145" = "dnasr281@gmail.com, System"
    "G146" = "dnasr281@gmail.com, System"
    "G148" = "dnasr281@gmail.com, System"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

